# Update the two-digit/one-digit division worksheet numbers to match
# output generated at 4250d90.
#
# Each table cell holds a single run with text like "59÷6=".
# We replace the dividend/divisor pairs with the new values. To avoid
# any risk of a later replacement accidentally re-matching text that
# was just written by an earlier replacement (some new values equal
# other old values), we do this in two passes: first swap every old
# value for a unique placeholder token, then swap every placeholder
# for its final new value.

$d = $word.ActiveDocument

# Phase 1: replace each original value with a unique placeholder token
# to avoid any chained collisions between old/new values.
$d.Content.Find.Execute("59÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER0@@", 2)
$d.Content.Find.Execute("46÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER1@@", 2)
$d.Content.Find.Execute("75÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER2@@", 2)
$d.Content.Find.Execute("43÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER3@@", 2)
$d.Content.Find.Execute("93÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER4@@", 2)
$d.Content.Find.Execute("71÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER5@@", 2)
$d.Content.Find.Execute("58÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER6@@", 2)
$d.Content.Find.Execute("96÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER7@@", 2)
$d.Content.Find.Execute("13÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER8@@", 2)
$d.Content.Find.Execute("95÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER9@@", 2)
$d.Content.Find.Execute("33÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER10@@", 2)
$d.Content.Find.Execute("83÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER11@@", 2)
$d.Content.Find.Execute("16÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER12@@", 2)
$d.Content.Find.Execute("36÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER13@@", 2)
$d.Content.Find.Execute("59÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER14@@", 2)
$d.Content.Find.Execute("43÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER15@@", 2)
$d.Content.Find.Execute("25÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER16@@", 2)
$d.Content.Find.Execute("87÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER17@@", 2)
$d.Content.Find.Execute("72÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER18@@", 2)
$d.Content.Find.Execute("92÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER19@@", 2)
$d.Content.Find.Execute("32÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER20@@", 2)
$d.Content.Find.Execute("48÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER21@@", 2)
$d.Content.Find.Execute("23÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER22@@", 2)
$d.Content.Find.Execute("30÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER23@@", 2)
$d.Content.Find.Execute("84÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER24@@", 2)

# Phase 2: replace each placeholder token with the final new value
$d.Content.Find.Execute("@@PLACEHOLDER0@@", $true, $false, $false, $false, $false, $true, 1, $false, "44÷3=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER1@@", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER2@@", $true, $false, $false, $false, $false, $true, 1, $false, "22÷3=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER3@@", $true, $false, $false, $false, $false, $true, 1, $false, "28÷5=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER4@@", $true, $false, $false, $false, $false, $true, 1, $false, "55÷5=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER5@@", $true, $false, $false, $false, $false, $true, 1, $false, "45÷5=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER6@@", $true, $false, $false, $false, $false, $true, 1, $false, "66÷6=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER7@@", $true, $false, $false, $false, $false, $true, 1, $false, "72÷7=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER8@@", $true, $false, $false, $false, $false, $true, 1, $false, "95÷7=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER9@@", $true, $false, $false, $false, $false, $true, 1, $false, "51÷4=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER10@@", $true, $false, $false, $false, $false, $true, 1, $false, "29÷4=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER11@@", $true, $false, $false, $false, $false, $true, 1, $false, "83÷6=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER12@@", $true, $false, $false, $false, $false, $true, 1, $false, "70÷7=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER13@@", $true, $false, $false, $false, $false, $true, 1, $false, "81÷5=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER14@@", $true, $false, $false, $false, $false, $true, 1, $false, "11÷2=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER15@@", $true, $false, $false, $false, $false, $true, 1, $false, "90÷4=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER16@@", $true, $false, $false, $false, $false, $true, 1, $false, "29÷6=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER17@@", $true, $false, $false, $false, $false, $true, 1, $false, "16÷6=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER18@@", $true, $false, $false, $false, $false, $true, 1, $false, "46÷8=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER19@@", $true, $false, $false, $false, $false, $true, 1, $false, "20÷2=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER20@@", $true, $false, $false, $false, $false, $true, 1, $false, "90÷6=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER21@@", $true, $false, $false, $false, $false, $true, 1, $false, "28÷8=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER22@@", $true, $false, $false, $false, $false, $true, 1, $false, "43÷4=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER23@@", $true, $false, $false, $false, $false, $true, 1, $false, "94÷2=", 2)
$d.Content.Find.Execute("@@PLACEHOLDER24@@", $true, $false, $false, $false, $false, $true, 1, $false, "89÷8=", 2)
